$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.051644320458384511
$ws.Range("B1").Value = -0.05164432047223172
$ws.Range("A2").Value = -0.012642373872628131
$ws.Range("B2").Value = 0.012642373860786543
$ws.Range("A3").Value = -0.02119337501181344
$ws.Range("B3").Value = 0.021193374990534
$ws.Range("A4").Value = -0.04922999833962801
$ws.Range("B4").Value = 0.049229998309049498
